# The underlying data rows for this observation export got re-sorted.
# Rows 4-8 keep all of their "static" columns (C, I, K, P, S, T, U, V, W,
# Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY, ...) in place, but the
# per-observation columns (A, B, D, E, F, G, H, Q, R) are permuted across
# the five rows as follows (new row -> source row):
#   4 <- 5
#   5 <- 4
#   6 <- 7
#   7 <- 8
#   8 <- 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually change between the rows.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

$firstRow = 4
$lastRow = 8

# new row number -> old (source) row number
$mapping = @{ 4 = 5; 5 = 4; 6 = 7; 7 = 8; 8 = 6 }

# Snapshot the current ("before") values for the relevant columns/rows
# first, since we are going to overwrite cells in place and some of the
# source rows depend on each other (e.g. 4 <-> 5).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowValues
}

# Now write the permuted values back.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $mapping[$r]
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = $srcValues[$col]
    }
}
